# Nalco aluminium ingot price sheet: a new circular (28-08-2025) was published,
# so it becomes the new first data row (Sl.no. 4) and the three previous
# entries shift down by one row (their Sl.no. values decrement accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlinks; they will be recreated in the correct
# row order below (this avoids stale refs left behind by a row insert).
$ws.Hyperlinks.Delete()

# Shift the three existing data rows (2-4) down to (3-5), bottom-up so we
# don't clobber a row before it has been copied.
$ws.Range("A4:F4").Copy($ws.Range("A5:F5"))
$ws.Range("A3:F3").Copy($ws.Range("A4:F4"))
$ws.Range("A2:F2").Copy($ws.Range("A3:F3"))

# Write the new latest-circular row into row 2.
$ws.Cells.Item(2, 1).Value = 4
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 271.05
$ws.Cells.Item(2, 5).Value = "28-08-2025"
$ws.Cells.Item(2, 6).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"

# Recreate the hyperlinks, in row order, against the refreshed link text.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")

# Hyperlinks.Add re-styles its target with the built-in "Hyperlink" look
# (underline + theme color); restore the plain centered look the rest of
# the table uses by re-pulling formatting from a same-row, already-styled
# cell.
$ws.Range("E2").Copy()
$ws.Range("F2:F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
